# TC32_Verify_store_location.xlsx — "Logic change for Logged in User"
#
# The login sequence used to be:
#   row3  CLICK      LoginOption
#   row4  ENTERTEXT  Uname1
#   row5  ENTERTEXT  Password1
#   row6  CLICK      LoginButton1
#
# It becomes:
#   row3  CLICK      LoginOption
#   row4  CLICK      LoginURL      <- new step
#   row5  ENTERTEXT  Uname         <- renamed object (was Uname1)
#   row6  ENTERTEXT  Password      <- renamed object (was Password1)
#   row7  CLICK      LoginButton   <- renamed object (was LoginButton1)
#
# Every following row simply shifts down by one (rows 7-16 -> 8-17) with
# their values/styles carried along automatically by the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- insert the new step as row 4, pushing the old rows 4-16 to 5-17 ---
$ws.Rows.Item(4).Insert()

# the freshly inserted row has no formatting yet - give it the same
# look as the row above it (border box, D/E bold like the other new
# login-url / empty cells in this table)
$ws.Range("A4:E4").Borders.LineStyle = 1
$ws.Range("D4:E4").Font.Bold = $true

# --- new row 4: CLICK / LoginURL / CSS / (blank) ---
$ws.Cells.Item(4, 1).Value = ""
$ws.Cells.Item(4, 2).Value = "CLICK"
$ws.Cells.Item(4, 3).Value = "LoginURL"
$ws.Cells.Item(4, 4).Value = "CSS"
$ws.Cells.Item(4, 5).Value = ""

# --- rows that shifted down one slot: fix up the renamed Object values ---
$ws.Cells.Item(5, 3).Value = "Uname"        # was "Uname1"
$ws.Cells.Item(6, 3).Value = "Password"     # was "Password1"
$ws.Cells.Item(7, 3).Value = "LoginButton"  # was "LoginButton1"

# --- selection moves to the newly-relevant ENTERTEXT/CLICK object cells ---
$ws.Activate()
$ws.Range("C5:C7").Select()
